$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 3300   # was 3287
$ws.Range("F5").Value = 209   # was 206
$ws.Range("F6").Value = 4802   # was 4787
$ws.Range("F7").Value = 461   # was 457
$ws.Range("F8").Value = 290   # was 289
$ws.Range("F9").Value = 172   # was 169
$ws.Range("F10").Value = 623   # was 621
$ws.Range("F12").Value = 23   # was 20
$ws.Range("F14").Value = 649   # was 648
$ws.Range("F15").Value = 286   # was 285
$ws.Range("F16").Value = 22   # was 20
$ws.Range("F17").Value = 88   # was 87
$ws.Range("F18").Value = 141   # was 140
$ws.Range("F19").Value = 339   # was 336
$ws.Range("F20").Value = 4735   # was 4726
$ws.Range("F24").Value = 5876   # was 5867
$ws.Range("F26").Value = 1192   # was 1190
$ws.Range("F28").Value = 665   # was 660
$ws.Range("F29").Value = 4416   # was 4412
$ws.Range("F31").Value = 89   # was 87
$ws.Range("F32").Value = 122   # was 121
$ws.Range("F33").Value = 851   # was 841
$ws.Range("F34").Value = 70   # was 67
$ws.Range("F36").Value = 775   # was 764
$ws.Range("F37").Value = 814   # was 800

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 37   # was 36
$ws.Range("F4").Value = 12   # was 11

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 1088   # was 1087

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1088   # was 1087
$ws.Range("F8").Value = 3300   # was 3287
$ws.Range("F9").Value = 209   # was 206
$ws.Range("F10").Value = 4802   # was 4787
$ws.Range("F11").Value = 461   # was 457
$ws.Range("F12").Value = 290   # was 289
$ws.Range("F13").Value = 172   # was 169
$ws.Range("F14").Value = 623   # was 621
$ws.Range("F16").Value = 23   # was 20
$ws.Range("F18").Value = 649   # was 648
$ws.Range("F19").Value = 286   # was 285
$ws.Range("F20").Value = 22   # was 20
$ws.Range("F21").Value = 37   # was 36
$ws.Range("F22").Value = 88   # was 87
$ws.Range("F23").Value = 141   # was 140
$ws.Range("F24").Value = 339   # was 336
$ws.Range("F25").Value = 4735   # was 4726
$ws.Range("F29").Value = 5876   # was 5867
$ws.Range("F31").Value = 1192   # was 1190
$ws.Range("F33").Value = 665   # was 660
$ws.Range("F34").Value = 4416   # was 4412
$ws.Range("F36").Value = 12   # was 11
$ws.Range("F37").Value = 89   # was 87
$ws.Range("F38").Value = 122   # was 121
$ws.Range("F39").Value = 851   # was 841
$ws.Range("F40").Value = 70   # was 67
$ws.Range("F42").Value = 775   # was 764
$ws.Range("F43").Value = 814   # was 800
